$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 5108.25
$ws.Range("J29").Value = 5145.0835
$ws.Range("L29").Value = 15435.2505
$ws.Range("N29").Value = -15997.2505
$ws.Range("H38").Value = 4294.3335
$ws.Range("J38").Value = 5087.8823
$ws.Range("L38").Value = 15263.6469
$ws.Range("N38").Value = -16007.6469
$ws.Range("H43").Value = 12728.571
$ws.Range("J43").Value = 13155.2
$ws.Range("L43").Value = 13155.2
$ws.Range("N43").Value = -13293.2
$ws.Range("H51").Value = 102177630
$ws.Range("I51").Value = 255436640
$ws.Range("J51").Value = 4966.6665
$ws.Range("K51").Value = 255436640
$ws.Range("L51").Value = 4966.6665
$ws.Range("M51").Value = -255436156
$ws.Range("N51").Value = -5934.6665
$ws.Range("H58").Value = 3416.3914
$ws.Range("J58").Value = 5428.357
$ws.Range("L58").Value = 16285.071
$ws.Range("N58").Value = -16585.071
$ws.Range("H112").Value = 669917.7
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 717697.5
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 2153092.5
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -2155308.5
$ws.Range("H135").Value = 2790.5715
$ws.Range("I135").Value = 1390.1
$ws.Range("K135").Value = 12510.9
$ws.Range("M135").Value = -9975.9
$ws.Range("H138").Value = 107276.33
$ws.Range("I138").Value = 2242.389
$ws.Range("J138").Value = 131208.11
$ws.Range("K138").Value = 6727.167
$ws.Range("L138").Value = 393624.33
$ws.Range("M138").Value = -1587.167
$ws.Range("N138").Value = -403904.33

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8821.965
$ws.Range("I32").Value = 8576.867
$ws.Range("J32").Value = 18993.5
$ws.Range("K32").Value = 8576.867
$ws.Range("L32").Value = 18993.5
$ws.Range("M32").Value = -8289.867
$ws.Range("N32").Value = -19567.5
$ws.Range("H74").Value = 2202.827
$ws.Range("I74").Value = 1350.7894
$ws.Range("K74").Value = 1350.7894
$ws.Range("M74").Value = -476.7893999999999
$ws.Range("H77").Value = 2202.827
$ws.Range("I77").Value = 1350.7894
$ws.Range("K77").Value = 6753.946999999999
$ws.Range("M77").Value = -2385.946999999999
$ws.Range("H132").Value = 2383.4036
$ws.Range("I132").Value = 2099.689
$ws.Range("K132").Value = 6299.066999999999
$ws.Range("M132").Value = -3769.066999999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").Value = $null

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2595.9614
$ws.Range("I31").Value = 2071.9524
$ws.Range("K31").Value = 2071.9524
$ws.Range("M31").Value = -1776.9524
$ws.Range("H33").Value = 35186.668
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").Value = $null
$ws.Range("H34").Value = 2595.9614
$ws.Range("I34").Value = 2071.9524
$ws.Range("K34").Value = 2071.9524
$ws.Range("M34").Value = -1869.9524

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 26580.47
$ws.Range("I68").Value = 62323.855
$ws.Range("J68").Value = 1560.1
$ws.Range("K68").Value = 186971.565
$ws.Range("L68").Value = 4680.299999999999
$ws.Range("M68").Value = -186160.565
$ws.Range("N68").Value = -6302.299999999999
$ws.Range("H71").Value = 26580.47
$ws.Range("I71").Value = 62323.855
$ws.Range("J71").Value = 1560.1
$ws.Range("K71").Value = 560914.6950000001
$ws.Range("L71").Value = 14040.9
$ws.Range("M71").Value = -556858.6950000001
$ws.Range("N71").Value = -22152.9
$ws.Range("H102").Value = 14278.286
$ws.Range("I102").Value = 3500
$ws.Range("K102").Value = 10500
$ws.Range("M102").Value = -8066
$ws.Range("H107").Value = 4850.421
$ws.Range("J107").Value = 5362.2354
$ws.Range("L107").Value = 16086.7062
$ws.Range("N107").Value = -19926.7062
$ws.Range("H113").Value = 2594.0588
$ws.Range("J113").Value = 2978.2856
$ws.Range("L113").Value = 8934.856800000001
$ws.Range("N113").Value = -13274.8568
$ws.Range("H122").Value = 995.26666
$ws.Range("I122").Value = 877.8889
$ws.Range("K122").Value = 7901.0001
$ws.Range("M122").Value = -5451.0001
$ws.Range("H129").Value = 84713.836
$ws.Range("I129").Value = 100776.8
$ws.Range("K129").Value = 302330.4
$ws.Range("M129").Value = -297330.4
$ws.Range("H131").Value = 133342.28
$ws.Range("I131").Value = 1112581.1
$ws.Range("J131").Value = 1802.7313
$ws.Range("K131").Value = 3337743.3
$ws.Range("L131").Value = 5408.1939
$ws.Range("M131").Value = -3332703.3
$ws.Range("N131").Value = -15488.1939

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 32221.5
$ws.Range("I26").Value = 29999
$ws.Range("K26").Value = 29999
$ws.Range("M26").Value = -29719
$ws.Range("H50").Value = 32221.5
$ws.Range("I50").Value = 29999
$ws.Range("K50").Value = 29999
$ws.Range("M50").Value = -29501
$ws.Range("H80").Value = 4007
$ws.Range("J80").Value = 4041.5833
$ws.Range("L80").Value = 4041.5833
$ws.Range("N80").Value = -6037.5833
$ws.Range("H83").Value = 4007
$ws.Range("J83").Value = 4041.5833
$ws.Range("L83").Value = 20207.9165
$ws.Range("N83").Value = -30191.9165
$ws.Range("H132").Value = 4436.5
$ws.Range("I132").Value = 4679.7
$ws.Range("J132").Value = 2004.5
$ws.Range("K132").Value = 14039.1
$ws.Range("L132").Value = 6013.5
$ws.Range("M132").Value = -11509.1
$ws.Range("N132").Value = -11073.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 28000
$ws.Range("J29").Value = 28000
$ws.Range("L29").Value = 28000
$ws.Range("N29").Value = -28590
$ws.Range("H33").Value = 18670
$ws.Range("I33").Value = 2000
$ws.Range("J33").Value = 27005
$ws.Range("K33").Value = 2000
$ws.Range("L33").Value = 27005
$ws.Range("M33").Value = -1710
$ws.Range("N33").Value = -27585
$ws.Range("H82").Value = 7083.8945
$ws.Range("I82").Value = 20387.666
$ws.Range("J82").Value = 943.6923
$ws.Range("K82").Value = 20387.666
$ws.Range("L82").Value = 943.6923
$ws.Range("M82").Value = -20026.666
$ws.Range("N82").Value = -1665.6923
$ws.Range("H85").Value = 7083.8945
$ws.Range("I85").Value = 20387.666
$ws.Range("J85").Value = 943.6923
$ws.Range("K85").Value = 20387.666
$ws.Range("L85").Value = 943.6923
$ws.Range("M85").Value = -19139.666
$ws.Range("N85").Value = -3439.6923
$ws.Range("H108").Value = 80748.5
$ws.Range("J108").Value = 73499.5
$ws.Range("L108").Value = 73499.5
$ws.Range("N108").Value = -81179.5
$ws.Range("H111").Value = 50000
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").Value = $null
$ws.Range("H132").Value = 3619.0579
$ws.Range("J132").Value = 3754.6206
$ws.Range("L132").Value = 11263.8618
$ws.Range("N132").Value = -16323.8618

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = $null
$ws.Range("H75").Value = 71250
$ws.Range("I75").Value = 55000
$ws.Range("K75").Value = 55000
$ws.Range("M75").Value = -54064
$ws.Range("H78").Value = 71250
$ws.Range("I78").Value = 55000
$ws.Range("K78").Value = 165000
$ws.Range("M78").Value = -160320
$ws.Range("H136").Value = 4539.967
$ws.Range("I136").Value = 5335.273
$ws.Range("J136").Value = 2352.875
$ws.Range("K136").Value = 16005.819
$ws.Range("L136").Value = 7058.625
$ws.Range("M136").Value = -13455.819
$ws.Range("N136").Value = -12158.625
